$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.907.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4657"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3651"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07362"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.380"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07126"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.503"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008690"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.919.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.045.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.131"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.255"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08907"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7546"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.474"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.909"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.177"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5269"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.326"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1653"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.427"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4846"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.09"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.657"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06289"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "
